$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.043.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.567.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.68%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.572.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.96%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.162"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.23%  "
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.020.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "59.086.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.12%  "
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.573.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.16%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "335.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +7.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.42%  "
$ws.Range("E26").Value = "  +0.54%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0767"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "157.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "18.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.890"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.851"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.21%  "
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "290.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.98%  "
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("E45").Value = "  -1.48%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.589"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0528"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("E49").Value = "  -1.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.939.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.13%  "
